$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values scraped on Tue Feb 14 21:43:57 UTC 2023

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "296.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.29%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.68%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.041"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.21%"

# Row 5
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.75%"

# Row 6
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.383"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.30%"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.596"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.26%"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9271"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.78%"

# Row 9
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.421"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.63%"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1208"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.20%"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1838"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.41%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08975"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.83%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03942"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-5.68%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1049"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.26%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001282"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.43%"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005918"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.56%"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.335"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.87%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3320"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.34%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.933"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.66%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1420"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.80%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3001"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.02%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04055"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.86%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001265"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.37%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003987"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "4.82%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001230"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.02%"

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.10%"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02408"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.05%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05209"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.03%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006402"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-4.31%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007801"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.59%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1327"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.51%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007563"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.35%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007837"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.78%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3214"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "10.89%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006787"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.80%"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.13%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04632"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "117.61%"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004203"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.13%"

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.13%"
